$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the bottom of the table.
# Column A holds date labels stored as text (same style as the existing
# "dd-mm-yyyy" text rows above), so we force a text number format before
# assigning the value to stop Excel from auto-converting the string into a
# date serial number, then restore the default ("Normal") style so the
# cells keep the same (unstyled) look as the rest of the table.

$ws.Range("A148").NumberFormat = "@"
$ws.Range("A148").Value = "03-08-2021"
$ws.Range("A148").Style = "Normal"
$ws.Range("B148").Value = 3317
$ws.Range("C148").Value = 7412

$ws.Range("A149").NumberFormat = "@"
$ws.Range("A149").Value = "04-08-2021"
$ws.Range("A149").Style = "Normal"
$ws.Range("B149").Value = 3463
$ws.Range("C149").Value = 8145
